# Apply the "revisions" edit to TableS1_stock_selection.xlsx
#
# Summary of the change:
#   - The note "Removed 2 population w/ wild dynamics" (a near-duplicate,
#     unused-alone shared string) is dropped.
#   - The note "Removed 2 population with wild dynamics" is expanded into a
#     full explanatory sentence, and now appears both on the "All" sheet
#     (row 29) and on the "Non-fish predator" sheet (row 9).
#   - The active sheet/selection moves: "All" becomes the selected tab with
#     A29 selected, and "Non-fish predator" is left with A10 selected (no
#     longer the active tab).

$wb = $excel.ActiveWorkbook

$newNote = "Removed 2 populations preventing model convergence because they exhibited population dynamics wildly divergent from stationary logistic population growth"

$wsAll = $wb.Worksheets.Item("All")
$wsAll.Range("A29").Value = $newNote

$wsNonFish = $wb.Worksheets.Item("Non-fish predator")
$wsNonFish.Range("A9").Value = $newNote

# Update selection on "Non-fish predator" first (it loses the active-tab
# flag once we activate "All" afterwards).
$wsNonFish.Activate()
$wsNonFish.Range("A10").Select()

# "All" ends up as the active/selected tab with A29 selected.
$wsAll.Activate()
$wsAll.Range("A29").Select()
